$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")

# Helper: a cell known to carry the default (unstyled) format for these data rows
$defaultStyleCell = $ws.Range("A123")

function Set-TextValue($cell, $value) {
    # Force the cell to remain plain text (avoid Excel auto-converting
    # date-looking strings like "2026-02-08" into date serial numbers).
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $defaultStyleCell.Style
}

# Row 123 - S122 Feed Additive ROI
Set-TextValue $ws.Range("B123") "Feed Additive ROI 3:1"
Set-TextValue $ws.Range("C123") "Article Summary"
Set-TextValue $ws.Range("D123") "sources/articles/S122_wattagnet_phytogenic_roi.pdf"
Set-TextValue $ws.Range("E123") "https://www.wattagnet.com/animal-health/article/15535560/phytogenic-feed-additives-deliver-roi"
Set-TextValue $ws.Range("F123") "2026-02-08"
Set-TextValue $ws.Range("G123") "2026-02-08"
Set-TextValue $ws.Range("H123") "Agent"
Set-TextValue $ws.Range("I123") "Primary web capture PDF. Supporting captures: sources/articles/S122_feedandadditive_phytogenic_roi.pdf; sources/articles/S122_ew_nutrition_phytogenics.pdf. Legacy summary retained at sources/articles/FeedAdditive_ROI_3to1.txt."

# Row 124 - S123 Urban vs Suburban Pet Habits
Set-TextValue $ws.Range("B124") "Urban vs Suburban Pet Habits"
Set-TextValue $ws.Range("C124") "Article Summary"
Set-TextValue $ws.Range("D124") "sources/articles/S123_petfoodindustry_urban_suburban.pdf"
Set-TextValue $ws.Range("E124") "https://www.petfoodindustry.com/nutrition/article/15468763/urban-vs-suburban-pet-owners-purchasing-habits"
Set-TextValue $ws.Range("F124") "2026-02-08"
Set-TextValue $ws.Range("G124") "2026-02-08"
Set-TextValue $ws.Range("H124") "Agent"
Set-TextValue $ws.Range("I124") "Primary web capture PDF. Legacy summary retained at sources/articles/PetFoodInd_UrbanSuburban.txt."

# Row 125 - S124 MARA Announcement 194 Summary
Set-TextValue $ws.Range("B125") "MARA Announcement 194 Summary"
Set-TextValue $ws.Range("C125") "Regulatory Summary"
Set-TextValue $ws.Range("D125") "sources/regulatory/S124_moa_announcement_194.pdf"
Set-TextValue $ws.Range("E125") "http://www.moa.gov.cn/govpublic/xmsyj/201912/t20191227_6334005.htm"
Set-TextValue $ws.Range("F125") "2020-07-01"
Set-TextValue $ws.Range("G125") "2026-02-08"
Set-TextValue $ws.Range("H125") "Agent"
Set-TextValue $ws.Range("I125") "Primary regulatory web capture PDF. Supporting captures: sources/regulatory/S124_feedstrategy_china_agp_ban.pdf; sources/regulatory/S124_mordor_china_feed_additives_market.pdf. Legacy summary retained at sources/regulatory/MARA_Announcement_194_Summary.txt."

# Row 126 - S125 Sector Deal Multiples 2020-2024
Set-TextValue $ws.Range("B126") "Sector Deal Multiples 2020-2024"
Set-TextValue $ws.Range("C126") "Transaction Summary"
Set-TextValue $ws.Range("D126") "sources/reports/S125_prnewswire_hh_zesty_paws.pdf"
Set-TextValue $ws.Range("E126") "https://www.prnewswire.com/news-releases/hh-group-acquires-zesty-paws-301361559.html"
Set-TextValue $ws.Range("F126") "2026-02-08"
Set-TextValue $ws.Range("G126") "2026-02-08"
Set-TextValue $ws.Range("H126") "Agent"
Set-TextValue $ws.Range("I126") "Primary transaction web capture PDF. Supporting captures: sources/reports/S125_generalmills_bluebuffalo.pdf; sources/reports/S125_swedencare_naturvet_press.pdf; sources/reports/S125_zoetis_mfa_phibro.pdf; sources/reports/S125_dsm_erber_group.pdf. Legacy summary retained at sources/reports/Sector_Deal_Multiples_2020-2024.txt. Historical Nasdaq link used in prior summary is no longer accessible."

# Row 127 - S126 EU Green Claims Directive Summary
Set-TextValue $ws.Range("B127") "EU Green Claims Directive Summary"
Set-TextValue $ws.Range("C127") "Regulatory Summary"
Set-TextValue $ws.Range("D127") "sources/regulatory/S126_ec_green_claims.pdf"
Set-TextValue $ws.Range("E127") "https://environment.ec.europa.eu/topics/circular-economy/green-claims_en"
Set-TextValue $ws.Range("F127") "2023-03-22"
Set-TextValue $ws.Range("G127") "2026-02-08"
Set-TextValue $ws.Range("H127") "Agent"
Set-TextValue $ws.Range("I127") "Primary EU web capture PDF. Supporting capture: sources/regulatory/S126_europarl_green_claims_train.pdf. Legacy summary retained at sources/regulatory/EU_Green_Claims_Directive_Summary.txt."

# Row 128 - S127 Nutrigenomics Review Summary
Set-TextValue $ws.Range("B128") "Nutrigenomics Review Summary"
Set-TextValue $ws.Range("C128") "Academic Summary"
Set-TextValue $ws.Range("D128") "sources/academic/S127_frontiers_nutrigenomics_review.pdf"
Set-TextValue $ws.Range("E128") "https://www.frontiersin.org/journals/veterinary-science/articles/10.3389/fvets.2020.00346/full"
Set-TextValue $ws.Range("F128") "2026-02-08"
Set-TextValue $ws.Range("G128") "2026-02-08"
Set-TextValue $ws.Range("H128") "Agent"
Set-TextValue $ws.Range("I128") "Primary academic web capture PDF. Supporting capture: sources/academic/S127_ncbi_pmc7575754.pdf. Legacy summary retained at sources/academic/Nutrigenomics_Review_Summary.txt."
